$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = 10
$ws.Range("A2").Value = "Test"
$ws.Range("C2").Value = 1
[void]$ws.Range("C2").Select()
